# Natmi following Dr Hou advice
# Update Ebi3-Il6st LR-pair sheet: recompute rows 2-5 for Sending cluster "ECs"
# and append rows 6-9 for Sending cluster "M2" (same Ligand/Receptor, all Target clusters).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ebi3"
$ws.Cells.Item(2, 3).Value = "Il6st"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1351216666666667
$ws.Cells.Item(2, 8).Value = 0.405365
$ws.Cells.Item(2, 9).Value = 0.01567245347423418
$ws.Cells.Item(2, 10).Value = 0.01567245347423418
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 44.017783
$ws.Cells.Item(2, 14).Value = 132.053349
$ws.Cells.Item(2, 15).Value = 0.253483683026081
$ws.Cells.Item(2, 16).Value = 0.253483683026081
$ws.Cells.Item(2, 17).Value = 5.947756201931667
$ws.Cells.Item(2, 18).Value = 53.529805817385
$ws.Cells.Item(2, 19).Value = 0.003972711228703779
$ws.Cells.Item(2, 20).Value = 0.003972711228703778

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ebi3"
$ws.Cells.Item(3, 3).Value = "Il6st"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1351216666666667
$ws.Cells.Item(3, 8).Value = 0.405365
$ws.Cells.Item(3, 9).Value = 0.01567245347423418
$ws.Cells.Item(3, 10).Value = 0.01567245347423418
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 99.76728333333334
$ws.Cells.Item(3, 14).Value = 299.30185
$ws.Cells.Item(3, 15).Value = 0.5745264004968147
$ws.Cells.Item(3, 16).Value = 0.5745264004968147
$ws.Cells.Item(3, 17).Value = 13.48072160280556
$ws.Cells.Item(3, 18).Value = 121.32649442525
$ws.Cells.Item(3, 19).Value = 0.009004238281505561
$ws.Cells.Item(3, 20).Value = 0.009004238281505561

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ebi3"
$ws.Cells.Item(4, 3).Value = "Il6st"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1351216666666667
$ws.Cells.Item(4, 8).Value = 0.405365
$ws.Cells.Item(4, 9).Value = 0.01567245347423418
$ws.Cells.Item(4, 10).Value = 0.01567245347423418
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.852095
$ws.Cells.Item(4, 14).Value = 32.556285
$ws.Cells.Item(4, 15).Value = 0.06249358376701795
$ws.Cells.Item(4, 16).Value = 0.06249358376701795
$ws.Cells.Item(4, 17).Value = 1.466353163225
$ws.Cells.Item(4, 18).Value = 13.197178469025
$ws.Cells.Item(4, 19).Value = 0.0009794277840267452
$ws.Cells.Item(4, 20).Value = 0.0009794277840267452

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ebi3"
$ws.Cells.Item(5, 3).Value = "Il6st"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1351216666666667
$ws.Cells.Item(5, 8).Value = 0.405365
$ws.Cells.Item(5, 9).Value = 0.01567245347423418
$ws.Cells.Item(5, 10).Value = 0.01567245347423418
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 19.01418566666667
$ws.Cells.Item(5, 14).Value = 57.042557
$ws.Cells.Item(5, 15).Value = 0.1094963327100864
$ws.Cells.Item(5, 16).Value = 0.1094963327100864
$ws.Cells.Item(5, 17).Value = 2.569228457589444
$ws.Cells.Item(5, 18).Value = 23.123056118305
$ws.Cells.Item(5, 19).Value = 0.001716076179998096
$ws.Cells.Item(5, 20).Value = 0.001716076179998096

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Ebi3"
$ws.Cells.Item(6, 3).Value = "Il6st"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.486480999999999
$ws.Cells.Item(6, 8).Value = 25.459443
$ws.Cells.Item(6, 9).Value = 0.9843275465257658
$ws.Cells.Item(6, 10).Value = 0.9843275465257658
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 44.017783
$ws.Cells.Item(6, 14).Value = 132.053349
$ws.Cells.Item(6, 15).Value = 0.253483683026081
$ws.Cells.Item(6, 16).Value = 0.253483683026081
$ws.Cells.Item(6, 17).Value = 373.556079091623
$ws.Cells.Item(6, 18).Value = 3362.004711824607
$ws.Cells.Item(6, 19).Value = 0.2495109717973773
$ws.Cells.Item(6, 20).Value = 0.2495109717973772

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Ebi3"
$ws.Cells.Item(7, 3).Value = "Il6st"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.486480999999999
$ws.Cells.Item(7, 8).Value = 25.459443
$ws.Cells.Item(7, 9).Value = 0.9843275465257658
$ws.Cells.Item(7, 10).Value = 0.9843275465257658
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 99.76728333333334
$ws.Cells.Item(7, 14).Value = 299.30185
$ws.Cells.Item(7, 15).Value = 0.5745264004968147
$ws.Cells.Item(7, 16).Value = 0.5745264004968147
$ws.Cells.Item(7, 17).Value = 846.67315442995
$ws.Cells.Item(7, 18).Value = 7620.05838986955
$ws.Cells.Item(7, 19).Value = 0.5655221622153092
$ws.Cells.Item(7, 20).Value = 0.5655221622153092

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Ebi3"
$ws.Cells.Item(8, 3).Value = "Il6st"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.486480999999999
$ws.Cells.Item(8, 8).Value = 25.459443
$ws.Cells.Item(8, 9).Value = 0.9843275465257658
$ws.Cells.Item(8, 10).Value = 0.9843275465257658
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.852095
$ws.Cells.Item(8, 14).Value = 32.556285
$ws.Cells.Item(8, 15).Value = 0.06249358376701795
$ws.Cells.Item(8, 16).Value = 0.06249358376701795
$ws.Cells.Item(8, 17).Value = 92.09609802769499
$ws.Cells.Item(8, 18).Value = 828.8648822492551
$ws.Cells.Item(8, 19).Value = 0.0615141559829912
$ws.Cells.Item(8, 20).Value = 0.06151415598299121

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Ebi3"
$ws.Cells.Item(9, 3).Value = "Il6st"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.486480999999999
$ws.Cells.Item(9, 8).Value = 25.459443
$ws.Cells.Item(9, 9).Value = 0.9843275465257658
$ws.Cells.Item(9, 10).Value = 0.9843275465257658
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 19.01418566666667
$ws.Cells.Item(9, 14).Value = 57.042557
$ws.Cells.Item(9, 15).Value = 0.1094963327100864
$ws.Cells.Item(9, 16).Value = 0.1094963327100864
$ws.Cells.Item(9, 17).Value = 161.363525390639
$ws.Cells.Item(9, 18).Value = 1452.271728515751
$ws.Cells.Item(9, 19).Value = 0.1077802565300883
$ws.Cells.Item(9, 20).Value = 0.1077802565300883

